# Apply weekly update to Fruta / Terminal La Palmera de La Serena - Membrillo sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44309
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 305000
$ws.Range("O4").Value = 310000
$ws.Range("P4").Value = 307500
$ws.Range("R4").Value = "Provincia de Cachapoal"
$ws.Range("S4").Value = 683

# Row 5
$ws.Range("D5").Value = 44309
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 285000
$ws.Range("O5").Value = 290000
$ws.Range("P5").Value = 287500
$ws.Range("R5").Value = "Provincia de Cachapoal"
$ws.Range("S5").Value = 639

# Row 6
$ws.Range("L6").Value = "Segunda"
$ws.Range("N6").Value = 255000
$ws.Range("O6").Value = 260000
$ws.Range("P6").Value = 257500
$ws.Range("S6").Value = 572

# Row 7
$ws.Range("D7").Value = 44295
$ws.Range("N7").Value = 225000
$ws.Range("O7").Value = 230000
$ws.Range("P7").Value = 227500
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 506

# Row 8
$ws.Range("D8").Value = 44295
$ws.Range("M8").Value = 16
$ws.Range("N8").Value = 195000
$ws.Range("O8").Value = 200000
$ws.Range("P8").Value = 197500
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 439

# Row 9
$ws.Range("D9").Value = 44316
$ws.Range("L9").Value = "Especial"
$ws.Range("N9").Value = 255000
$ws.Range("O9").Value = 260000
$ws.Range("P9").Value = 257500
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 572

# Row 10
$ws.Range("D10").Value = 44316
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 225000
$ws.Range("O10").Value = 230000
$ws.Range("P10").Value = 227500
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 506

# Row 11
$ws.Range("D11").Value = 44294
$ws.Range("L11").Value = "Primera"
$ws.Range("N11").Value = 225000
$ws.Range("O11").Value = 230000
$ws.Range("P11").Value = 227500
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 506

# Row 12
$ws.Range("D12").Value = 44294
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 16
$ws.Range("N12").Value = 195000
$ws.Range("O12").Value = 200000
$ws.Range("P12").Value = 197500
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 439
